# Edit map_test.xlsx: add two new "luck" (运气) event rows to the map sheet,
# mirroring the existing "新闻" (news) event rows, and fix the sheet's
# frozen-pane / selection metadata.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert "运气1" row right before the current row 16 (old "杭州" row) ---
$ws.Rows.Item(16).Insert()
$ws.Range("A16").Value = 15
$ws.Range("B16").Value = "运气1"
$ws.Range("M16").Value = "抽取1张运气卡。"

# --- Insert "运气2" row right before the current row 33 (old "澳门" row,
#     which after the first insertion above now sits at row 33) ---
$ws.Rows.Item(33).Insert()
$ws.Range("A33").Value = 32
$ws.Range("B33").Value = "运气2"
$ws.Range("M33").Value = "抽取1张运气卡。"

# --- Update the view state (frozen pane / active selection) ---
$ws.Application.ActiveWindow.ScrollRow = 22
$ws.Range("J41").Select()
